# This script applies updated TPM-derived metric values to the
# "LR-pairs_lrc2p" sheet (Icam1-Itgal), columns G:T, rows 2-26.
# The new values come from a refreshed TPM computation (see commit
# message "update scripts wuth new tpm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab-free CSV block holding the new values for range G2:T26
# (columns: Ligand avg expr, Ligand total expr, Ligand derived spec (avg),
#  Ligand derived spec (total), Receptor-expressing cells, Receptor
#  detection rate, Receptor avg expr, Receptor total expr, Receptor
#  derived spec (avg), Receptor derived spec (total), Edge avg expr
#  weight, Edge total expr weight, Edge avg expr derived spec,
#  Edge total expr derived spec)
$csv = @"
35.160799,105.482397,0.2238945559395223,0.2238945559395223,3.0,1.0,0.3547066666666667,1.06412,0.006934574868045491,0.00693457486804549,12.47176981062667,112.24592829564,0.001552613560710416,0.001552613560710416
35.160799,105.482397,0.2238945559395223,0.2238945559395223,3.0,1.0,0.519977,1.559931,0.01016563762403213,0.01016563762403213,18.282806781623,164.545261034607,0.002276030921674774,0.002276030921674774
35.160799,105.482397,0.2238945559395223,0.2238945559395223,3.0,1.0,26.50170333333334,79.50511,0.5181127482682332,0.5181127482682331,931.8210640609633,8386.38957654867,0.1160026237001215,0.1160026237001215
35.160799,105.482397,0.2238945559395223,0.2238945559395223,3.0,1.0,0.09443866666666667,0.283316,0.001846291784117559,0.001846291784117558,3.320538976494666,29.884850788452,0.000413374679139789,0.000413374679139789
35.160799,105.482397,0.2238945559395223,0.2238945559395223,3.0,1.0,23.67963033333333,71.03889099999999,0.4629407474555717,0.4629407474555717,832.5947225446362,7493.352502901726,0.1036499130778758,0.1036499130778757
32.208719,96.626157,0.2050965007332699,0.2050965007332699,3.0,1.0,0.3547066666666667,1.06412,0.006934574868045491,0.00693457486804549,11.42464735409333,102.82182618684,0.001422257039509007,0.001422257039509007
32.208719,96.626157,0.2050965007332699,0.2050965007332699,3.0,1.0,0.519977,1.559931,0.01016563762403213,0.01016563762403213,16.747793079463,150.730137715167,0.002084936704411462,0.002084936704411461
32.208719,96.626157,0.2050965007332699,0.2050965007332699,3.0,1.0,26.50170333333334,79.50511,0.5181127482682332,0.5181127482682331,853.5859156846968,7682.273241162271,0.1062631116551122,0.1062631116551121
32.208719,96.626157,0.2050965007332699,0.2050965007332699,3.0,1.0,0.09443866666666667,0.283316,0.001846291784117559,0.001846291784117558,3.041748477401334,27.375736296612,0.000378667984255097,0.0003786679842550969
32.208719,96.626157,0.2050965007332699,0.2050965007332699,3.0,1.0,23.67963033333333,71.03889099999999,0.4629407474555717,0.4629407474555717,762.6905594302096,6864.215034871887,0.09494752734998217,0.09494752734998216
51.53356533333334,154.600696,0.3281519491717758,0.3281519491717758,3.0,1.0,0.3547066666666667,1.06412,0.006934574868045491,0.00693457486804549,18.27929918083556,164.51369262752,0.002275594259626738,0.002275594259626738
51.53356533333334,154.600696,0.3281519491717758,0.3281519491717758,3.0,1.0,0.519977,1.559931,0.01016563762403213,0.01016563762403213,26.79626870133067,241.166418311976,0.003335873800900084,0.003335873800900083
51.53356533333334,154.600696,0.3281519491717758,0.3281519491717758,3.0,1.0,26.50170333333334,79.50511,0.5181127482682332,0.5181127482682331,1365.727260172951,12291.54534155656,0.1700197082349663,0.1700197082349663
51.53356533333334,154.600696,0.3281519491717758,0.3281519491717758,3.0,1.0,0.09443866666666667,0.283316,0.001846291784117559,0.001846291784117558,4.866761198659556,43.800850787936,0.0006058642476980124,0.0006058642476980123
51.53356533333334,154.600696,0.3281519491717758,0.3281519491717758,3.0,1.0,23.67963033333333,71.03889099999999,0.4629407474555717,0.4629407474555717,1220.295776852015,10982.66199166814,0.1519149086285847,0.1519149086285847
0.5955593333333333,1.786678,0.003792362411113143,0.003792362411113143,3.0,1.0,0.3547066666666667,1.06412,0.006934574868045491,0.00693457486804549,0.2112488659288889,1.90123979336,0.00002629842106662561,0.0000262984210666256
0.5955593333333333,1.786678,0.003792362411113143,0.003792362411113143,3.0,1.0,0.519977,1.559931,0.01016563762403213,0.01016563762403213,0.3096771554686667,2.787094399218,0.00003855178201037698,0.00003855178201037697
0.5955593333333333,1.786678,0.003792362411113143,0.003792362411113143,3.0,1.0,26.50170333333334,79.50511,0.5181127482682332,0.5181127482682331,15.78333676939778,142.05003092458,0.001964871311250974,0.001964871311250973
0.5955593333333333,1.786678,0.003792362411113143,0.003792362411113143,3.0,1.0,0.09443866666666667,0.283316,0.001846291784117559,0.001846291784117558,0.05624382936088889,0.506194464248,0.000007001807562034451,0.00000700180756203445
0.5955593333333333,1.786678,0.003792362411113143,0.003792362411113143,3.0,1.0,23.67963033333333,71.03889099999999,0.4629407474555717,0.4629407474555717,14.10262485489978,126.923623694098,0.001755639089223133,0.001755639089223133
37.54313466666667,112.629404,0.2390646317443189,0.2390646317443189,3.0,1.0,0.3547066666666667,1.06412,0.006934574868045491,0.00693457486804549,13.31680015383111,119.85120138448,0.001657811587132704,0.001657811587132704
37.54313466666667,112.629404,0.2390646317443189,0.2390646317443189,3.0,1.0,0.519977,1.559931,0.01016563762403213,0.01016563762403213,19.52156653456933,175.694098811124,0.002430244415035435,0.002430244415035435
37.54313466666667,112.629404,0.2390646317443189,0.2390646317443189,3.0,1.0,26.50170333333334,79.50511,0.5181127482682332,0.5181127482682331,994.9570171393823,8954.613154254439,0.1238624333667822,0.1238624333667822
37.54313466666667,112.629404,0.2390646317443189,0.2390646317443189,3.0,1.0,0.09443866666666667,0.283316,0.001846291784117559,0.001846291784117558,3.545523580407111,31.909712223664,0.0004413830654626257,0.0004413830654626257
37.54313466666667,112.629404,0.2390646317443189,0.2390646317443189,3.0,1.0,23.67963033333333,71.03889099999999,0.4629407474555717,0.4629407474555717,889.0075504612182,8001.067954150963,0.110672759309906,0.110672759309906
"@

$lines = $csv -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$rowCount = $lines.Count
$colCount = ($lines[0] -split ",").Count

$values = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $fields = $lines[$i] -split ","
    for ($j = 0; $j -lt $colCount; $j++) {
        $values[$i, $j] = [double]$fields[$j]
    }
}

$ws.Range("G2:T$(1 + $rowCount)").Value = $values

Write-Host "Updated $rowCount rows x $colCount cols in range G2:T$(1 + $rowCount)"
